$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R1").Value = 1103
$ws.Range("S1").Value = 1160

$ws.Range("R2").Value = 226
$ws.Range("S2").Value = 244

$ws.Range("R3").Value = 888
$ws.Range("S3").Value = 943

$ws.Range("R4").Value = 534
$ws.Range("S4").Value = 578

$ws.Range("R5").Value = 213
$ws.Range("S5").Value = 244

$ws.Range("R6").Value = 55
$ws.Range("S6").Value = 58

$ws.Range("R7").Value = 412
$ws.Range("S7").Value = 490

$ws.Range("R1:S7").Font.Color = 0

$ws.Application.ActiveWindow.ScrollColumn = 15
$ws.Range("S7").Select()
$ws.Range("A1:S7").Select()

$wb.Windows.Item(1).left = $wb.Windows.Item(1).left
